$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8295.15
$ws.Range("I116").Value = 10991.917
$ws.Range("J116").Value = 4250
$ws.Range("K116").Value = 10991.917
$ws.Range("L116").Value = 4250
$ws.Range("M116").Value = -7549.916999999999
$ws.Range("N116").Value = -11134

$ws.Range("H129").Value = 731.2553
$ws.Range("J129").Value = 958
$ws.Range("L129").Value = 2874
$ws.Range("N129").Value = -12874

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1031.0834
$ws.Range("I122").Value = 809.8
$ws.Range("J122").Value = 1189.1428
$ws.Range("K122").Value = 2429.4
$ws.Range("L122").Value = 3567.4284
$ws.Range("M122").Value = 20.60000000000036
$ws.Range("N122").Value = -8467.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 916.6667
$ws.Range("I107").Value = 666.6667
$ws.Range("K107").Value = 666.6667
$ws.Range("M107").Value = 1253.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2923.025
$ws.Range("I31").Value = 1386.92
$ws.Range("K31").Value = 1386.92
$ws.Range("M31").Value = -1091.92

$ws.Range("H34").Value = 2923.025
$ws.Range("I34").Value = 1386.92
$ws.Range("K34").Value = 1386.92
$ws.Range("M34").Value = -1184.92

$ws.Range("H99").Value = 2462.6
$ws.Range("I99").Value = 2387.5557
$ws.Range("J99").Value = 2618.4614
$ws.Range("K99").Value = 2387.5557
$ws.Range("L99").Value = 2618.4614
$ws.Range("M99").Value = -889.5556999999999
$ws.Range("N99").Value = -5614.4614

$ws.Range("H122").Value = 1095.619
$ws.Range("I122").Value = 997.46155
$ws.Range("J122").Value = 1255.125
$ws.Range("K122").Value = 2992.38465
$ws.Range("L122").Value = 3765.375
$ws.Range("M122").Value = -542.38465
$ws.Range("N122").Value = -8665.375

$ws.Range("H126").Value = 2462.6
$ws.Range("I126").Value = 2387.5557
$ws.Range("J126").Value = 2618.4614
$ws.Range("K126").Value = 7162.6671
$ws.Range("L126").Value = 7855.3842
$ws.Range("M126").Value = -4692.6671
$ws.Range("N126").Value = -12795.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3832.4517
$ws.Range("J5").Value = 13387.625
$ws.Range("L5").Value = 40162.875
$ws.Range("N5").Value = -40386.875

$ws.Range("H39").Value = 2500
$ws.Range("J39").Value = 2500
$ws.Range("L39").Value = 7500
$ws.Range("N39").Value = -8088

$ws.Range("H58").Value = 1830
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 2245
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 6735
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -6991

$ws.Range("H64").Value = 1897048.5
$ws.Range("I64").Value = 1404
$ws.Range("J64").Value = 3371438.8
$ws.Range("K64").Value = 4212
$ws.Range("L64").Value = 10114316.4
$ws.Range("M64").Value = -3942
$ws.Range("N64").Value = -10114856.4

$ws.Range("H67").Value = 1897048.5
$ws.Range("I67").Value = 1404
$ws.Range("J67").Value = 3371438.8
$ws.Range("K67").Value = 4212
$ws.Range("L67").Value = 10114316.4
$ws.Range("M67").Value = -3276
$ws.Range("N67").Value = -10116188.4

$ws.Range("H82").Value = 62502260
$ws.Range("I82").Value = 787.75
$ws.Range("J82").Value = 83336080
$ws.Range("K82").Value = 2363.25
$ws.Range("L82").Value = 250008240
$ws.Range("M82").Value = -1957.25
$ws.Range("N82").Value = -250009052

$ws.Range("H85").Value = 62502260
$ws.Range("I85").Value = 787.75
$ws.Range("J85").Value = 83336080
$ws.Range("K85").Value = 2363.25
$ws.Range("L85").Value = 250008240
$ws.Range("M85").Value = -959.25
$ws.Range("N85").Value = -250011048

$ws.Range("H88").Value = 3277.8125
$ws.Range("J88").Value = 3277.8125
$ws.Range("L88").Value = 9833.4375
$ws.Range("N88").Value = -10689.4375

$ws.Range("H91").Value = 3277.8125
$ws.Range("J91").Value = 3277.8125
$ws.Range("L91").Value = 9833.4375
$ws.Range("N91").Value = -12797.4375

$ws.Range("H94").Value = 2287.6
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 2609.5
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 7828.5
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -9180.5

$ws.Range("H97").Value = 6935.5557
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 7777.5
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 23332.5
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -24324.5

$ws.Range("H103").Value = 2679.8572
$ws.Range("I103").Value = 399.5
$ws.Range("J103").Value = 3592
$ws.Range("K103").Value = 1198.5
$ws.Range("L103").Value = 10776
$ws.Range("M103").Value = -319.5
$ws.Range("N103").Value = -12534

$ws.Range("H106").Value = 4495.1763
$ws.Range("J106").Value = 4495.1763
$ws.Range("L106").Value = 13485.5289
$ws.Range("N106").Value = -15377.5289

$ws.Range("H109").Value = 1244.3334
$ws.Range("I109").Value = 454.7857
$ws.Range("J109").Value = 2094.6155
$ws.Range("K109").Value = 1364.3571
$ws.Range("L109").Value = 6283.8465
$ws.Range("M109").Value = -324.3571000000002
$ws.Range("N109").Value = -8363.8465

$ws.Range("H112").Value = 56207772
$ws.Range("I112").Value = 4900
$ws.Range("J112").Value = 59513824
$ws.Range("K112").Value = 14700
$ws.Range("L112").Value = 178541472
$ws.Range("M112").Value = -13592
$ws.Range("N112").Value = -178543688

$ws.Range("H122").Value = 726.8461
$ws.Range("I122").Value = 512.5
$ws.Range("J122").Value = 822.1111
$ws.Range("K122").Value = 4612.5
$ws.Range("L122").Value = 7398.9999
$ws.Range("M122").Value = -2162.5
$ws.Range("N122").Value = -12298.9999

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0

$ws.Range("H131").Value = 81589280
$ws.Range("I131").Value = 520
$ws.Range("J131").Value = 103840760
$ws.Range("K131").Value = 1560
$ws.Range("L131").Value = 311522280
$ws.Range("M131").Value = 3480
$ws.Range("N131").Value = -311532360

$ws.Range("H135").Value = 3832.4517
$ws.Range("J135").Value = 13387.625
$ws.Range("L135").Value = 120488.625
$ws.Range("N135").Value = -125558.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4579.8887
$ws.Range("I70").Value = 4705.9
$ws.Range("J70").Value = 4219.857
$ws.Range("K70").Value = 4705.9
$ws.Range("L70").Value = 4219.857
$ws.Range("M70").Value = -4435.9
$ws.Range("N70").Value = -4759.857

$ws.Range("H73").Value = 4579.8887
$ws.Range("I73").Value = 4705.9
$ws.Range("J73").Value = 4219.857
$ws.Range("K73").Value = 4705.9
$ws.Range("L73").Value = 4219.857
$ws.Range("M73").Value = -3769.9
$ws.Range("N73").Value = -6091.857

$ws.Range("H122").Value = 1040.3043
$ws.Range("I122").Value = 923.5
$ws.Range("J122").Value = 1222
$ws.Range("K122").Value = 2770.5
$ws.Range("L122").Value = 3666
$ws.Range("M122").Value = -320.5
$ws.Range("N122").Value = -8566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2285.1155
$ws.Range("I7").Value = 1724.2941
$ws.Range("J7").Value = 3344.4443
$ws.Range("K7").Value = 1724.2941
$ws.Range("L7").Value = 3344.4443
$ws.Range("M7").Value = -1612.2941
$ws.Range("N7").Value = -3568.4443

$ws.Range("H22").Value = 1209.95
$ws.Range("I22").Value = 657.3333
$ws.Range("J22").Value = 1446.7858
$ws.Range("K22").Value = 657.3333
$ws.Range("L22").Value = 1446.7858
$ws.Range("M22").Value = -362.3333
$ws.Range("N22").Value = -2036.7858

$ws.Range("H27").Value = 1209.95
$ws.Range("I27").Value = 657.3333
$ws.Range("J27").Value = 1446.7858
$ws.Range("K27").Value = 657.3333
$ws.Range("L27").Value = 1446.7858
$ws.Range("M27").Value = -550.3333
$ws.Range("N27").Value = -1660.7858

$ws.Range("H126").Value = 2285.1155
$ws.Range("I126").Value = 1724.2941
$ws.Range("J126").Value = 3344.4443
$ws.Range("K126").Value = 5172.8823
$ws.Range("L126").Value = 10033.3329
$ws.Range("M126").Value = -2702.8823
$ws.Range("N126").Value = -14973.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4313.8
$ws.Range("I122").Value = 2536.3635
$ws.Range("J122").Value = 5710.357
$ws.Range("K122").Value = 7609.0905
$ws.Range("L122").Value = 17131.071
$ws.Range("M122").Value = -5159.0905
$ws.Range("N122").Value = -22031.071
